$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = -1
$ws.Range("F2").Value = 8
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 9
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 2
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 34
$ws.Range("E7").Value = 14
$ws.Range("F7").Value = 11
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 20
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 11
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = -2
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 4
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 1
$ws.Range("E14").Value = 11
$ws.Range("F14").Value = 15
$ws.Range("E15").Value = -1
$ws.Range("F15").Value = 7
$ws.Range("E16").Value = -1
$ws.Range("F16").Value = 2
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = -5
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = -5
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 36
$ws.Range("E22").Value = 14
$ws.Range("F22").Value = -3
$ws.Range("E23").Value = 14
$ws.Range("F23").Value = -3
$ws.Range("E24").Value = 6
$ws.Range("F24").Value = 20
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 11
$ws.Range("E26").Value = 6
$ws.Range("F26").Value = 20
$ws.Range("E27").Value = 6
$ws.Range("F27").Value = 15
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 7
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = 36
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = 5
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 3
$ws.Range("E32").Value = 6
$ws.Range("F32").Value = 20
$ws.Range("E33").Value = 12
$ws.Range("F33").Value = 4
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 0
$ws.Range("E35").Value = 3
$ws.Range("F35").Value = 3
$ws.Range("E36").Value = 7
$ws.Range("F36").Value = 4
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = 5
$ws.Range("E38").Value = 6
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = 34
$ws.Range("E40").Value = 12
$ws.Range("F40").Value = 7
$ws.Range("E41").Value = 14
$ws.Range("F41").Value = -3
$ws.Range("E42").Value = 6
$ws.Range("F42").Value = 5
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 9
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 1
$ws.Range("F45").Value = 25
$ws.Range("E46").Value = 6
$ws.Range("F46").Value = 6
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 6
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = 19
$ws.Range("E49").Value = 9
$ws.Range("F49").Value = 34
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = 27
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 9
$ws.Range("F52").Value = 3
$ws.Range("E53").Value = 3
$ws.Range("F53").Value = 3
$ws.Range("E54").Value = 3
$ws.Range("F54").Value = 3
$ws.Range("E55").Value = 3
$ws.Range("F55").Value = 0
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 3
$ws.Range("E57").Value = 3
$ws.Range("F57").Value = 3
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 3
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 7
$ws.Range("E60").Value = 5
$ws.Range("F60").Value = 5
$ws.Range("E61").Value = 11
$ws.Range("F61").Value = 0
$ws.Range("E62").Value = 12
$ws.Range("F62").Value = 0
$ws.Range("E63").Value = 12
$ws.Range("F63").Value = 17
$ws.Range("E64").Value = 17
$ws.Range("F64").Value = -6
$ws.Range("E65").Value = 11
$ws.Range("F65").Value = 12
$ws.Range("E66").Value = 8
$ws.Range("F66").Value = 1
$ws.Range("E67").Value = 7
$ws.Range("F67").Value = 22
$ws.Range("E68").Value = 15
$ws.Range("F68").Value = -4
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 6
$ws.Range("E70").Value = 21
$ws.Range("F70").Value = -2
$ws.Range("E71").Value = 14
$ws.Range("F71").Value = 10
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = 8
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = 4
$ws.Range("E74").Value = 17
$ws.Range("F74").Value = 36
$ws.Range("E75").Value = 9
$ws.Range("F75").Value = 2
$ws.Range("E76").Value = 13
$ws.Range("F76").Value = 71
$ws.Range("F77").Value = 36
$ws.Range("E78").Value = 37
$ws.Range("F78").Value = 14
$ws.Range("E79").Value = 1
$ws.Range("F79").Value = 5
$ws.Range("E80").Value = 15
$ws.Range("F80").Value = 0
$ws.Range("E81").Value = 3
$ws.Range("F81").Value = 2
$ws.Range("E82").Value = 1
$ws.Range("F82").Value = 4
$ws.Range("E83").Value = 0
$ws.Range("F83").Value = 5
$ws.Range("E84").Value = 2
$ws.Range("F84").Value = 1
$ws.Range("E85").Value = 7
$ws.Range("F85").Value = 34
$ws.Range("E86").Value = 2
$ws.Range("F86").Value = 1
$ws.Range("E87").Value = 11
$ws.Range("F87").Value = 0
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 3
$ws.Range("E89").Value = 12
$ws.Range("F89").Value = 16
$ws.Range("E90").Value = 19
$ws.Range("F90").Value = -8
$ws.Range("E91").Value = -1
$ws.Range("F91").Value = 5
$ws.Range("E92").Value = 3
$ws.Range("F92").Value = 4
$ws.Range("E93").Value = 13
$ws.Range("F93").Value = -2
$ws.Range("E94").Value = 17
$ws.Range("F94").Value = 36
$ws.Range("E95").Value = 6
$ws.Range("F95").Value = 20
$ws.Range("E96").Value = -8
$ws.Range("F96").Value = 12
$ws.Range("E97").Value = 1
$ws.Range("F97").Value = 3
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 4
$ws.Range("E99").Value = 10
$ws.Range("F99").Value = 2
$ws.Range("E100").Value = 37
$ws.Range("F100").Value = 14
$ws.Range("E101").Value = 6
$ws.Range("F101").Value = -2
$ws.Range("E102").Value = 3
$ws.Range("F102").Value = 3
$ws.Range("E103").Value = -10
$ws.Range("F103").Value = 23
$ws.Range("E104").Value = -10
$ws.Range("F104").Value = 23
$ws.Range("E105").Value = 14
$ws.Range("F105").Value = 33
$ws.Range("E106").Value = 14
$ws.Range("F106").Value = 6
$ws.Range("E107").Value = 7
$ws.Range("F107").Value = 37
$ws.Range("E108").Value = 0
$ws.Range("F108").Value = -6
$ws.Range("E109").Value = 3
$ws.Range("F109").Value = 0
$ws.Range("F110").Value = 3
$ws.Range("E111").Value = 15
$ws.Range("F111").Value = 16
$ws.Range("E112").Value = 0
$ws.Range("F112").Value = 5
$ws.Range("E113").Value = 16
$ws.Range("F113").Value = 20
$ws.Range("E114").Value = 15
$ws.Range("F114").Value = -4
$ws.Range("E115").Value = 5
$ws.Range("F115").Value = 26
$ws.Range("E116").Value = 0
$ws.Range("F116").Value = 29
$ws.Range("E117").Value = 10
$ws.Range("F117").Value = 1
$ws.Range("E118").Value = 13
$ws.Range("F118").Value = 49
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = 21
$ws.Range("E120").Value = 23
$ws.Range("F120").Value = 9
$ws.Range("E121").Value = 23
$ws.Range("F121").Value = 9
$ws.Range("E122").Value = 9
$ws.Range("F122").Value = -1
$ws.Range("E123").Value = 19
$ws.Range("F123").Value = -8
$ws.Range("E124").Value = 6
$ws.Range("F124").Value = 5
$ws.Range("E125").Value = 5
$ws.Range("F125").Value = 6
$ws.Range("E126").Value = 17
$ws.Range("F126").Value = -6
$ws.Range("E127").Value = 11
$ws.Range("F127").Value = 12
$ws.Range("E128").Value = 8
$ws.Range("F128").Value = 1
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 22
$ws.Range("E130").Value = -5
$ws.Range("F130").Value = 36
$ws.Range("E131").Value = -4
$ws.Range("F131").Value = 12
$ws.Range("E132").Value = 13
$ws.Range("F132").Value = -2
$ws.Range("E133").Value = 4
$ws.Range("F133").Value = 17
$ws.Range("E134").Value = 13
$ws.Range("F134").Value = 34
$ws.Range("E135").Value = 6
$ws.Range("F135").Value = 5
$ws.Range("E136").Value = 9
$ws.Range("F136").Value = 34
$ws.Range("E137").Value = 3
$ws.Range("F137").Value = 17
$ws.Range("E138").Value = 13
$ws.Range("F138").Value = 28
$ws.Range("E139").Value = 13
$ws.Range("F139").Value = 2
$ws.Range("E140").Value = 17
$ws.Range("F140").Value = 36
$ws.Range("E141").Value = 4
$ws.Range("F141").Value = 7
$ws.Range("E142").Value = 9
$ws.Range("F142").Value = 36
$ws.Range("E143").Value = 10
$ws.Range("F143").Value = 14
$ws.Range("E144").Value = 9
$ws.Range("F144").Value = 34
